# Adds a new price-snapshot column "2026-02-01 15:14:29" right before the
# existing "nom"/"url_produit" columns (DD/DE), shifting them one column to
# the right (DD->DE, DE->DF). For each product row, the new column is seeded
# with a copy of the most recent prior snapshot value (column DC, the last
# existing snapshot column) when that snapshot held a price; rows whose
# snapshot history had already run dry (no DC value) are left blank, same
# as the existing empty snapshot cells in that row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 206
$dcCol = 107   # column DC – last existing snapshot column before the edit
$newHeader = "2026-02-01 15:14:29"

# Capture the current (pre-insert) values of column DC for every data row,
# so we know what to carry forward into the freshly inserted column.
$dcVals = @{}
for ($r = 2; $r -le $lastRow; $r++) {
    $v = $ws.Cells.Item($r, $dcCol).Value2
    if ($v -ne $null -and $v -ne "") {
        $dcVals[$r] = $v
    }
}

# Insert the new column at DD; "nom" (was DD) and "url_produit" (was DE)
# shift right to DE/DF automatically, carrying their values/styles along.
$ws.Columns("DD").Insert()

# New header cell for the inserted column.
$ws.Range("DD1").Value = $newHeader

# Populate the new DD column for every row that had a live price in the old
# DC column; leave the rest blank (matching the already-empty history cells).
for ($r = 2; $r -le $lastRow; $r++) {
    if ($dcVals.ContainsKey($r)) {
        $ws.Cells.Item($r, 108).Value = $dcVals[$r]
    }
}
